$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set new column widths for F and G
$ws.Columns.Item(6).ColumnWidth = 20
$ws.Columns.Item(7).ColumnWidth = 15

# Header row
$ws.Range("F1").Value = "Full-Time Duration"
$ws.Range("G1").Value = "Flex Duration"

# Data rows - Full-Time Duration (F) and Flex Duration (G)
$durations = @(
    @{Row=2;  F="4 місяці";  G="8 місяців"},
    @{Row=3;  F="3 місяці";  G="5 місяців"},
    @{Row=4;  F="4 місяці";  G="7 місяців"},
    @{Row=5;  F="3 місяці";  G="5 місяців"},
    @{Row=6;  F="3 місяці";  G="4 місяці"},
    @{Row=7;  F="5 місяців"; G="7 місяців"},
    @{Row=8;  F="5 місяців"; G="8 місяців"},
    @{Row=9;  F="3 місяці";  G="4 місяці"},
    @{Row=10; F="3 місяці";  G="8 місяців"},
    @{Row=11; F="3 місяці";  G="4 місяці"}
)

foreach ($d in $durations) {
    $ws.Cells.Item($d.Row, 6).Value = $d.F
    $ws.Cells.Item($d.Row, 7).Value = $d.G
}

# Apply header style (same as other header cells) to F1:G1
$ws.Range("F1:G1").Style = $ws.Range("E1").Style

# Apply data row style (same as other data cells) to F2:G11
$ws.Range("F2:G11").Style = $ws.Range("E2").Style
